$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new data row (previously non-existent) is inserted at row 34; all rows
# that used to be 34..134 shift down to become 35..135. Insert a blank row at
# position 34 which pushes everything below it (including formatting) down
# by one, exactly reproducing that shift.
$ws.Rows.Item(34).Insert()

# Populate the freshly inserted row 34 with the new record. Every field
# except D (Fecha), K/L/M (prices) and P (Precio $/Kg) is identical to the
# record that used to occupy row 34 (now row 35), so reuse those values by
# copying from row 35 first and then overwriting the changed fields.
$ws.Cells.Item(34, 1).Value2 = $ws.Cells.Item(35, 1).Value2()
$ws.Cells.Item(34, 2).Value2 = $ws.Cells.Item(35, 2).Value2()
$ws.Cells.Item(34, 3).Value2 = $ws.Cells.Item(35, 3).Value2()
$ws.Cells.Item(34, 4).Value2 = 44497
$ws.Cells.Item(34, 5).Value2 = $ws.Cells.Item(35, 5).Value2()
$ws.Cells.Item(34, 6).Value2 = $ws.Cells.Item(35, 6).Value2()
$ws.Cells.Item(34, 7).Value2 = $ws.Cells.Item(35, 7).Value2()
$ws.Cells.Item(34, 8).Value2 = $ws.Cells.Item(35, 8).Value2()
$ws.Cells.Item(34, 9).Value2 = $ws.Cells.Item(35, 9).Value2()
$ws.Cells.Item(34, 10).Value2 = $ws.Cells.Item(35, 10).Value2()
$ws.Cells.Item(34, 11).Value2 = 7000
$ws.Cells.Item(34, 12).Value2 = 7000
$ws.Cells.Item(34, 13).Value2 = 7000
$ws.Cells.Item(34, 14).Value2 = $ws.Cells.Item(35, 14).Value2()
$ws.Cells.Item(34, 15).Value2 = $ws.Cells.Item(35, 15).Value2()
$ws.Cells.Item(34, 16).Value2 = 1167
$ws.Cells.Item(34, 17).Value2 = $ws.Cells.Item(35, 17).Value2()
$ws.Cells.Item(34, 18).Value2 = $ws.Cells.Item(35, 18).Value2()

$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
